# Repull data, push all data, mean calculation
# Update the "dSF" (column F) values for the rows whose data was refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 6
    3  = 3
    4  = 1
    6  = -4
    8  = 12
    9  = -2
    10 = -3
    11 = 2
    12 = 3
    13 = 8
    14 = 1
    15 = 4
    16 = -2
    17 = 1
    18 = 5
    19 = 7
    20 = -2
    22 = -1
    23 = -1
    24 = -1
    25 = 3
    27 = -3
    28 = 1
    29 = -1
    30 = 4
    31 = -4
    32 = -2
    33 = 1
    35 = 2
    36 = -2
    38 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
